# Refresh the crypto symbol list (Price / Volume(1h) columns, plus the
# BitKan/HotbitToken ranking swap in rows 24-25) to match the latest scrape.
# Target cells are formatted as Text ("@") before the write so values like
# "306.34" / "0.94%" are stored as literal strings (matching the original
# inlineStr cells) instead of being auto-coerced into numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.94%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.60%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.116"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.21%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07904"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.81%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.119"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.99%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.968"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.53%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9199"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.34%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09727"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.46%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1856"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.81%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08730"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.55%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03570"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.17%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09931"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.27%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001450"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.64%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005708"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.90%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.467"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.10%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.121"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.21%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.632"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "16.74%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3391"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.41%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.24%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.170"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.58%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2020"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-12.12%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04562"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.95%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.33%"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.005037"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.17%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "14.24%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01859"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.93%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04766"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.68%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007603"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.34%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1400"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.53%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007867"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.65%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002204"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.31%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.28%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006332"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.48%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.17%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.00%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.06"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "549.12%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.17%"
